$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 12 (the lone "어떤 데이터가 전달 될지 설계" note row) - all the
# rows below it shift up by one, matching the API section header moving
# from row 13 -> row 12, etc.
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(11).RowHeight = 17.25

# Fix remaining typo: "댓글의 id" -> "댓글 id" (now on row 25 after the shift)
$ws.Cells.Item(25, 5).Value = "댓글 id"

# Update the active selection to match the authored edit
$ws.Range("B21").Select()
